# edit.ps1
# Applies the "actualizacion de url en pdf" change set:
#   1. "CV online:" / "Web:" block: add the light-gray theme color to every
#      run + the paragraph mark run properties, merge "Web"+": " into a
#      single run "Web: ", merge the URL run with the trailing "home" run
#      (new target: .../presentacion/home), and drop the now-empty
#      paragraph (the one that only carried an underline rPr) that used to
#      follow it.
#   2. Move the <w:lastRenderedPageBreak/> marker from the start of the
#      "2018 - 2020 (Titulado)" run to the start of the
#      "Resumen de Tecnologías" run.
#   3. In the "Logros" paragraph of the first job, join the two runs that
#      made up "...reduciendo tiempos" / "operativos en un 20%." into one
#      run, and move the following <w:br/> + <w:lastRenderedPageBreak/> so
#      that the break now directly precedes "Implementación de Pruebas...".

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) "CV online: ... cv" + "Web: ... home" (+ drop trailing empty para)
# ---------------------------------------------------------------------

$rCv = $d.Content
$rCv.Find.Execute("CV online:", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rCv.Expand(4) | Out-Null

$rWeb = $d.Range($rCv.End, $d.Content.End)
$rWeb.Find.Execute("Web:", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rWeb.Expand(4) | Out-Null

$rEmpty = $d.Range($rWeb.End, $rWeb.End)
$rEmpty.Expand(4) | Out-Null

$rBlock = $d.Range($rCv.Start, $rEmpty.End)

$xmlBlock = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">
<w:body>
<w:p w14:paraId="218B97A5" w14:textId="70731460" w:rsidR="00A679A9" w:rsidRDefault="00A679A9">
  <w:pPr>
    <w:rPr>
      <w:i/>
      <w:iCs/>
      <w:color w:val="A6A6A6" w:themeColor="background1" w:themeShade="A6"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:color w:val="A6A6A6" w:themeColor="background1" w:themeShade="A6"/>
    </w:rPr>
    <w:t xml:space="preserve">CV online: </w:t>
  </w:r>
  <w:r w:rsidR="009D5924" w:rsidRPr="00FF688F">
    <w:rPr>
      <w:i/>
      <w:iCs/>
      <w:color w:val="A6A6A6" w:themeColor="background1" w:themeShade="A6"/>
    </w:rPr>
    <w:t>https://franco-bernal.github.io/presentacion/</w:t>
  </w:r>
  <w:r w:rsidR="006F4CAC">
    <w:rPr>
      <w:i/>
      <w:iCs/>
      <w:color w:val="A6A6A6" w:themeColor="background1" w:themeShade="A6"/>
    </w:rPr>
    <w:t>cv</w:t>
  </w:r>
</w:p>
<w:p w14:paraId="31421B3B" w14:textId="68C0B5F9" w:rsidR="006C0FEA" w:rsidRPr="006C0FEA" w:rsidRDefault="006C0FEA" w:rsidP="006C0FEA">
  <w:pPr>
    <w:rPr>
      <w:i/>
      <w:iCs/>
      <w:color w:val="A6A6A6" w:themeColor="background1" w:themeShade="A6"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:color w:val="A6A6A6" w:themeColor="background1" w:themeShade="A6"/>
    </w:rPr>
    <w:t xml:space="preserve">Web: </w:t>
  </w:r>
  <w:r w:rsidRPr="00FF688F">
    <w:rPr>
      <w:i/>
      <w:iCs/>
      <w:color w:val="A6A6A6" w:themeColor="background1" w:themeShade="A6"/>
    </w:rPr>
    <w:t>https://franco-bernal.github.io/presentacion/home</w:t>
  </w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

$rBlock.InsertXML($xmlBlock)

# ---------------------------------------------------------------------
# 2) Move <w:lastRenderedPageBreak/> from "2018 - 2020 (Titulado)" to
#    "Resumen de Tecnologías"
# ---------------------------------------------------------------------

$rTitulado = $d.Content
$rTitulado.Find.Execute("2018 - 2020 (Titulado)", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rTitulado.Expand(4) | Out-Null

$rResumenTec = $d.Range($rTitulado.End, $d.Content.End)
$rResumenTec.Find.Execute("Resumen de Tecnologías", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rResumenTec.Expand(4) | Out-Null

$rPageBreak = $d.Range($rTitulado.Start, $rResumenTec.End)

$xmlPageBreak = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">
<w:body>
<w:p w14:paraId="2F98E9FF" w14:textId="77777777" w:rsidR="002706F8" w:rsidRDefault="00000000">
  <w:r>
    <w:t>2018 - 2020 (Titulado)</w:t>
  </w:r>
</w:p>
<w:p w14:paraId="18559CE7" w14:textId="77777777" w:rsidR="002706F8" w:rsidRDefault="00000000">
  <w:pPr>
    <w:pStyle w:val="Ttulo1"/>
  </w:pPr>
  <w:r>
    <w:lastRenderedPageBreak/>
    <w:t>Resumen de Tecnologías</w:t>
  </w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

$rPageBreak.InsertXML($xmlPageBreak)

# ---------------------------------------------------------------------
# 3) "Logros" paragraph (Fábrica de Calzados Gino): merge the split
#    sentence about the 20% improvement and relocate the break +
#    lastRenderedPageBreak in front of "Implementación de Pruebas..."
# ---------------------------------------------------------------------

$rLogros = $d.Content
$rLogros.Find.Execute("Desarrollo de scripts en Node.js", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rLogros.Expand(4) | Out-Null

$xmlLogros = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">
<w:body>
<w:p w14:paraId="13D86551" w14:textId="77777777" w:rsidR="002706F8" w:rsidRDefault="00000000">
  <w:r>
    <w:t>Mantención y actualización de la infraestructura de ecommerce, lo que incrementó el rendimiento en un 30%.</w:t>
  </w:r>
  <w:r>
    <w:br/>
    <w:t>Optimización del backend, mejorando los tiempos de respuesta en un 40% en sitios ecommerce con tráfico alto.</w:t>
  </w:r>
  <w:r>
    <w:br/>
    <w:t>Desarrollo de scripts en Node.js y Python para automatizar tareas, reduciendo tiempos operativos en un 20%.</w:t>
  </w:r>
  <w:r>
    <w:br/>
  </w:r>
  <w:r>
    <w:lastRenderedPageBreak/>
    <w:t>Implementación de Pruebas de Penetración que aseguraron la protección contra vulnerabilidades de inyección SQL y otros riesgos.</w:t>
  </w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

$rLogros.InsertXML($xmlLogros)
